# Fruta / hortaliza, semanal
# Insert a new weekly price record at row 89 (pushing the existing
# rows 89-117 down to 90-118) and populate it with the new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 89..117 down by one to make room for the new record.
$ws.Rows("89:89").Insert()

# Populate the newly inserted row 89 with this week's Papaya price data.
$ws.Cells.Item(89, 1).Value  = 10
$ws.Cells.Item(89, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(89, 3).Value  = "La Araucanía"
$ws.Cells.Item(89, 4).Value  = 45146
$ws.Cells.Item(89, 5).Value  = 9
$ws.Cells.Item(89, 6).Value  = "Fruta"
$ws.Cells.Item(89, 7).Value  = 100108
$ws.Cells.Item(89, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(89, 9).Value  = 100108004
$ws.Cells.Item(89, 10).Value = "Papaya"
$ws.Cells.Item(89, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(89, 12).Value = "Primera"
$ws.Cells.Item(89, 13).Value = 50
$ws.Cells.Item(89, 14).Value = 25000
$ws.Cells.Item(89, 15).Value = 25000
$ws.Cells.Item(89, 16).Value = 25000
$ws.Cells.Item(89, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(89, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(89, 19).Value = 2500
$ws.Cells.Item(89, 20).Value = 10
